$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.2300186582838494
$ws.Range("E2").Value = 7.101739543477936
$ws.Range("F2").Value = 0.001512150340307394
$ws.Range("D3").Value = 0.222991587245047
$ws.Range("E3").Value = 8.230594652065671
$ws.Range("F3").Value = 0.001311467503360575
$ws.Range("D4").Value = 0.2008539547731077
$ws.Range("E4").Value = 8.852077964825687
$ws.Range("F4").Value = 0.001430755153499783
$ws.Range("D5").Value = 0.2176383991115586
$ws.Range("E5").Value = 9.768461413004566
$ws.Range("F5").Value = 0.006708721006691201
$ws.Range("D6").Value = 0.3074448830751838
$ws.Range("E6").Value = 10.71882365957675
$ws.Range("F6").Value = 0.02787433025776194
$ws.Range("D7").Value = 0.2848692949272273
$ws.Range("E7").Value = 13.9672509546112
$ws.Range("F7").Value = 0.0689364079333957
$ws.Range("D8").Value = 0.2301982192820561
$ws.Range("E8").Value = 14.47461954519428
$ws.Range("F8").Value = 0.1194018693226749
$ws.Range("D9").Value = 0.2288596421186697
$ws.Range("E9").Value = 14.02428988554301
$ws.Range("F9").Value = 0.08521290653765921
$ws.Range("D10").Value = 0.1062080410469003
$ws.Range("E10").Value = 12.30332272342335
$ws.Range("F10").Value = 0.03305091918572096
$ws.Range("D11").Value = 0.1075393528102039
$ws.Range("E11").Value = 13.01385173874577
$ws.Range("F11").Value = 0.06133888890874527
$ws.Range("D12").Value = 0.1101310488324591
$ws.Range("E12").Value = 13.50735043447484
$ws.Range("F12").Value = 0.06156072412074207
$ws.Range("D13").Value = 0.09131094886240723
$ws.Range("E13").Value = 14.12979131407006
$ws.Range("F13").Value = 0.05475046199935095
$ws.Range("D14").Value = 0.2176669619571188
$ws.Range("E14").Value = 14.63398270590704
$ws.Range("F14").Value = 0.1696557266145625
$ws.Range("D15").Value = 0.2119249067092369
$ws.Range("E15").Value = 14.90843498316629
$ws.Range("F15").Value = 0.2422359224449915
$ws.Range("D16").Value = 0.223633627004179
$ws.Range("E16").Value = 15.35619884448591
$ws.Range("F16").Value = 0.3053763296264926
$ws.Range("D17").Value = 0.1810302514322489
$ws.Range("E17").Value = 15.94958979062786
$ws.Range("F17").Value = 0.3843046185059165
$ws.Range("D18").Value = 0.2136737591899481
$ws.Range("E18").Value = 16.23607054969736
$ws.Range("F18").Value = 0.4960666013335486
$ws.Range("D19").Value = 0.2332859465516377
$ws.Range("E19").Value = 16.40027282850033
$ws.Range("F19").Value = 0.7818666558254169
$ws.Range("D20").Value = 0.2104442662843391
$ws.Range("E20").Value = 16.48637291124163
$ws.Range("F20").Value = 0.7356941455486997
$ws.Range("D21").Value = 0.2109897773499699
$ws.Range("E21").Value = 16.51963533745865
$ws.Range("F21").Value = 0.7897474847617028
$ws.Range("D22").Value = 0.2126454706557688
$ws.Range("E22").Value = 16.77964558369623
$ws.Range("F22").Value = 0.8649352785530209
$ws.Range("D23").Value = 0.2780205030253442
$ws.Range("E23").Value = 17.83105717399513
$ws.Range("F23").Value = 1.144712176558246
